# SAV-700: Update charts importer test fixture
# - Ensure ID of charting date recorded program data element is constant
# - Ensure ID of all complex chart core questions is constant

$wb = $excel.ActiveWorkbook

$core = $wb.Worksheets.Item("Core")
$testChart = $wb.Worksheets.Item("Test Chart")
$metadata = $wb.Worksheets.Item("Metadata")

# 1) Move the "highlighted" cell format from Core!R2 to Test Chart!A2 before
#    we normalise Core's R column formats below.
$core.Range("R2").Copy()
$testChart.Range("A2").PasteSpecial(-4122)

# 2) Update the "code" column on the Core sheet so every complex-chart core
#    question id is constant (matches its name instead of a generated
#    testchartcorecodeN placeholder).
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartType"
$core.Range("A5").Value = "ComplexChartSubtype"

# 3) Ensure the ID of the charting date recorded program data element is
#    constant on the Test Chart sheet.
$testChart.Range("A2").Value = "PatientChartingDate"

# 4) Normalise the format of the Core sheet's visibilityStatus (R) column so
#    rows 2, 3 and 5 match row 4's style.
$core.Range("R4").Copy()
$core.Range("R2").PasteSpecial(-4122)
$core.Range("R3").PasteSpecial(-4122)
$core.Range("R5").PasteSpecial(-4122)

# 5) Column width tweaks that came along with the content changes above.
$metadata.Columns.Item(1).ColumnWidth = 23.63
$testChart.Columns.Item(1).ColumnWidth = 19.13
